$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Top sale table (rows 2-7) ---
$ws.Range("G2").Value = "SSO Terget"

$ws.Range("C3").Value = 31670
$ws.Range("D3").Value = 47373

$ws.Range("C4").Value = 32594
$ws.Range("D4").Value = 20516
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 3

$ws.Range("C5").Value = 25704
$ws.Range("D5").Value = 16756

$ws.Range("B6").Value = "Arman Hosen"
$ws.Range("C6").Value = 31242
$ws.Range("D6").Value = 14214
$ws.Range("E6").Value = 2

# --- Daily stock headers (update the report date) ---
$ws.Range("A1").Value = "Mangrove Communication   20.01.2025"
$ws.Range("A10").Value = "DAILY STOCK                         (20/01/2025) "

# --- Stock table rows 13-31 ---
$ws.Range("C13").Value = 109658

$ws.Range("C14").Value = 240152
$ws.Range("D14").Value = 121210
$ws.Range("E14").Value = 124675

$ws.Range("C18").ClearContents()

$ws.Range("C20").Value = 9500
$ws.Range("D20").Value = 4500
$ws.Range("E20").ClearContents()

$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()

$ws.Range("C22").Value = 2110
$ws.Range("D22").Value = 1880
$ws.Range("E22").ClearContents()

$ws.Range("D24").Value = 5

$ws.Range("C25").Value = 31

$ws.Range("C26").Value = 55
$ws.Range("D26").Value = 3

$ws.Range("C27").Value = 92
$ws.Range("D27").Value = 3

# --- Sheet view: scroll position and active cell ---
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("J28").Select()
